# Updated cryptos list on Thu Apr 18 19:45:29 UTC 2024 with GitHub Actions
#
# Updates the "Price" (column D) and "Volume(1h)" (column E) values for the
# crypto rows on Sheet1. Values are written as literal text (matching the
# original inlineStr cell type) - the leading apostrophe forces Excel to
# treat numeric-looking strings (e.g. "140.00", "0.999") as text instead of
# silently coercing them to numbers, and resetting the style back to
# "Normal" afterwards avoids leaving a stray quote-prefixed number format
# applied to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-DValue "D2" '63.288.25'
$ws.Range("E2").Value = '  +3.91%  '
Set-DValue "D3" '3.059.85'
$ws.Range("E3").Value = '  +2.30%  '
$ws.Range("E4").Value = '  -0.11%  '
Set-DValue "D5" '548.48'
$ws.Range("E5").Value = '  +2.75%  '
Set-DValue "D6" '140.00'
$ws.Range("E6").Value = '  +4.87%  '
$ws.Range("E7").Value = '  -0.14%  '
Set-DValue "D8" '3.052.96'
$ws.Range("E8").Value = '  +2.30%  '
$ws.Range("E9").Value = '  +1.02%  '
Set-DValue "D10" '6.39'
$ws.Range("E10").Value = '  +5.21%  '
Set-DValue "D11" '0.150'
$ws.Range("E11").Value = '  +1.72%  '
$ws.Range("E12").Value = '  +1.98%  '
$ws.Range("E13").Value = '  +2.70%  '
Set-DValue "D14" '34.76'
$ws.Range("E14").Value = '  +2.00%  '
Set-DValue "D15" '3.554.97'
$ws.Range("E15").Value = '  +2.11%  '
Set-DValue "D16" '63.266.84'
$ws.Range("E16").Value = '  +3.58%  '
Set-DValue "D17" '3.060.18'
$ws.Range("E17").Value = '  +1.91%  '
Set-DValue "D19" '6.73'
$ws.Range("E19").Value = '  +2.22%  '
Set-DValue "D20" '481.78'
$ws.Range("E20").Value = '  +4.09%  '
Set-DValue "D21" '13.64'
$ws.Range("E21").Value = '  +3.53%  '
Set-DValue "D22" '0.673'
$ws.Range("E22").Value = '  -0.03%  '
Set-DValue "D23" '7.22'
$ws.Range("E23").Value = '  +4.15%  '
Set-DValue "D24" '80.63'
$ws.Range("E24").Value = '  +1.96%  '
Set-DValue "D25" '12.54'
$ws.Range("E25").Value = '  +4.41%  '
Set-DValue "D26" '0.999'
$ws.Range("E26").Value = '  -0.02%  '
Set-DValue "D27" '2.74'
$ws.Range("E27").Value = '  +2.62%  '
Set-DValue "D28" '7.91'
$ws.Range("E28").Value = '  +0.96%  '
Set-DValue "D29" '1.98'
$ws.Range("E29").Value = '  +4.75%  '
Set-DValue "D30" '1.00'
$ws.Range("E30").Value = '  -0.19%  '
Set-DValue "D31" '26.01'
$ws.Range("E31").Value = '  +2.27%  '
$ws.Range("E32").Value = '  +1.10%  '
Set-DValue "D33" '2.44'
$ws.Range("E33").Value = '  +7.19%  '
Set-DValue "D34" '5.67'
$ws.Range("E34").Value = '  +3.94%  '
Set-DValue "D35" '55.46'
$ws.Range("E35").Value = '  +0.20%  '
Set-DValue "D36" '5.97'
$ws.Range("E36").Value = '  +1.83%  '
Set-DValue "D37" '468.31'
$ws.Range("E37").Value = '  +2.93%  '
Set-DValue "D38" '0.0820'
$ws.Range("E38").Value = '  +4.46%  '
Set-DValue "D39" '0.0396'
$ws.Range("E39").Value = '  +3.51%  '
Set-DValue "D40" '3.071.39'
$ws.Range("E40").Value = '  -4.32%  '
$ws.Range("E41").Value = '  +0.41%  '
Set-DValue "D42" '8.25'
$ws.Range("E42").Value = '  +1.41%  '
Set-DValue "D43" '2.58'
$ws.Range("E43").Value = '  +5.24%  '
Set-DValue "D44" '28.03'
$ws.Range("E44").Value = '  +1.34%  '
Set-DValue "D45" '0.254'
$ws.Range("E45").Value = '  +3.63%  '
$ws.Range("E46").Value = '  -0.11%  '
Set-DValue "D47" '2.05'
$ws.Range("E47").Value = '  +2.78%  '
Set-DValue "D48" '0.109'
$ws.Range("E48").Value = '  +1.42%  '
Set-DValue "D49" '116.86'
$ws.Range("E49").Value = '  -1.90%  '
Set-DValue "D50" '0.0₃0508'
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("E51").Value = '  +2.69%  '
